$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data block: rows 2..49 (4 years x 12 months), columns A..H.
# Each year's 12-row chunk is currently ordered Jan..Dec and needs to be
# reordered to Oct, Nov, Dec, Jan, Feb, ..., Sep (i.e. rotate the 12 rows
# so October..December move to the front of that year's block).

$firstRow = 2
$lastRow = 49
$numRows = $lastRow - $firstRow + 1
$numCols = 8
$rowsPerYear = 12

$rng = $ws.Range("A$($firstRow):H$($lastRow)")
$vals = $rng.Value2

# New-Object arrays are 0-based; $vals from Value2 is 1-based (COM style).
$newvals = New-Object 'object[,]' $numRows, $numCols

$numYears = $numRows / $rowsPerYear
for ($y = 0; $y -lt $numYears; $y++) {
    $yearBase = $y * $rowsPerYear
    for ($m = 0; $m -lt $rowsPerYear; $m++) {
        # source month-of-year index (0-based, 0=Jan .. 11=Dec) that should
        # land at destination position $m within this year's block.
        $srcMonthIdx = ($m + 9) % $rowsPerYear
        $srcRow1 = $yearBase + $srcMonthIdx + 1   # 1-based row into $vals
        $dstRow0 = $yearBase + $m                 # 0-based row into $newvals
        for ($c = 0; $c -lt $numCols; $c++) {
            $newvals[$dstRow0, $c] = $vals[$srcRow1, $c + 1]
        }
    }
}

$rng.Value2 = $newvals
